# Auto-generated edit script: updates crypto price/volume(1h) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'290.61"
$ws.Range("E2").Value = "'-3.87%"
$ws.Range("D3").Value = "'30.76"
$ws.Range("E3").Value = "'-5.95%"
$ws.Range("D4").Value = "'4.955"
$ws.Range("E4").Value = "'-0.25%"
$ws.Range("D5").Value = "'0.07209"
$ws.Range("E5").Value = "'-7.97%"
$ws.Range("D6").Value = "'1.779"
$ws.Range("E6").Value = "'-8.81%"
$ws.Range("E7").Value = "'-2.50%"
$ws.Range("D8").Value = "'3.722"
$ws.Range("E8").Value = "'-2.03%"
$ws.Range("D9").Value = "'0.8958"
$ws.Range("E9").Value = "'-3.38%"
$ws.Range("D10").Value = "'0.1650"
$ws.Range("E10").Value = "'-6.90%"
$ws.Range("D11").Value = "'0.07702"
$ws.Range("E11").Value = "'-1.69%"
$ws.Range("D12").Value = "'0.07997"
$ws.Range("E12").Value = "'-7.79%"
$ws.Range("D13").Value = "'0.03036"
$ws.Range("E13").Value = "'-3.35%"
$ws.Range("E14").Value = "'-0.20%"
$ws.Range("D15").Value = "'0.001499"
$ws.Range("E15").Value = "'-1.01%"
$ws.Range("D16").Value = "'0.005711"
$ws.Range("E16").Value = "'-0.22%"
$ws.Range("D17").Value = "'3.478"
$ws.Range("E17").Value = "'0.44%"
$ws.Range("D18").Value = "'2.083"
$ws.Range("E18").Value = "'-3.29%"
$ws.Range("E19").Value = "'-0.47%"
$ws.Range("D20").Value = "'0.1319"
$ws.Range("E20").Value = "'0.15%"
$ws.Range("D21").Value = "'4.031"
$ws.Range("E21").Value = "'-6.32%"
$ws.Range("D23").Value = "'0.04516"
$ws.Range("E23").Value = "'-1.45%"
$ws.Range("E24").Value = "'-0.81%"
$ws.Range("D25").Value = "'0.004012"
$ws.Range("E25").Value = "'-9.71%"
$ws.Range("D26").Value = "'0.0001249"
$ws.Range("E26").Value = "'-0.10%"
$ws.Range("D39").Value = "'0.01600"
$ws.Range("E39").Value = "'-6.33%"
$ws.Range("D40").Value = "'0.04391"
$ws.Range("E40").Value = "'-7.72%"
$ws.Range("D41").Value = "'0.007302"
$ws.Range("E41").Value = "'-5.33%"
$ws.Range("E42").Value = "'-3.50%"
$ws.Range("D43").Value = "'0.007707"
$ws.Range("D44").Value = "'0.002069"
$ws.Range("E44").Value = "'-11.66%"
$ws.Range("D45").Value = "'0.009204"
$ws.Range("E45").Value = "'-20.81%"
$ws.Range("D46").Value = "'0.00005914"
$ws.Range("E46").Value = "'-5.55%"
$ws.Range("E47").Value = "'-0.11%"
$ws.Range("D48").Value = "'2.247"
$ws.Range("E48").Value = "'173.92%"
$ws.Range("E49").Value = "'-3.24%"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.11%"
$ws.Range("E51").Value = "'-0.11%"
